$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("0.05803533643484116", "0.9879117012023926", "0.1935855448246002", "0.9605362415313721"),
    @("0.009382420219480991", "0.9983578324317932", "0.1803392469882965", "0.961767315864563"),
    @("0.004995067603886127", "0.998677134513855", "0.1929948031902313", "0.9580056071281433"),
    @("0.002701822202652693", "0.9993613958358765", "0.1193691864609718", "0.9706586599349976"),
    @("0.001904790173284709", "0.9994525909423828", "0.2169000804424286", "0.9582107663154602"),
    @("0.001922387978993356", "0.9995894432067871", "0.1365380734205246", "0.9636139869689941"),
    @("0.00130282889585942", "0.9997262954711914", "0.173736035823822", "0.9697011113166809"),
    @("0.00154472969006747", "0.9996122717857361", "0.1994035094976425", "0.9598522782325745"),
    @("0.001021661213599145", "0.9997262954711914", "0.1475272923707962", "0.9763354063034058"),
    @("0.0005883114645257592", "0.9998631477355957", "0.2034564912319183", "0.969222366809845"),
    @("0.001231090049259365", "0.9997491240501404", "0.1697392016649246", "0.9756514430046082"),
    @("0.0009885343024507165", "0.9997262954711914", "0.1065578684210777", "0.9757198691368103"),
    @("0.0003313153283670545", "0.9999087452888489", "0.140577420592308", "0.9766773581504822"),
    @("0.000410029839258641", "0.9997718930244446", "0.1408723145723343", "0.9741467833518982"),
    @("0.0004322770400904119", "0.9999087452888489", "0.1230722367763519", "0.9786608219146729"),
    @("0.0004415939329192042", "0.9998859763145447", "0.1304000616073608", "0.9767457842826843"),
    @("0.0001019667179207318", "0.9999544024467468", "0.1218756511807442", "0.979002833366394"),
    @("0.00071460020262748", "0.9998631477355957", "0.1543091684579849", "0.9795499444007874"),
    @("0.0007278452394530177", "0.9997718930244446", "0.06904315948486328", "0.9848847389221191"),
    @("0.00012238509953022", "0.999977171421051", "0.117917388677597", "0.979755163192749"),
    @("1.852984860306606E-05", "1", "0.1316826790571213", "0.9768141508102417"),
    @("0.0003433595411479473", "0.9998859763145447", "0.1590884625911713", "0.9757882356643677"),
    @("0.0002469796454533935", "0.9999315738677979", "0.1491259783506393", "0.9782504439353943"),
    @("0.0003464491746854037", "0.9999544024467468", "0.1737604886293411", "0.9751726984977722"),
    @("0.0002032486372627318", "0.9999087452888489", "0.1370669007301331", "0.9789344072341919"),
    @("0.0005540908896364272", "0.9999315738677979", "0.1023864299058914", "0.9790711998939514"),
    @("2.155990478058811E-05", "1", "0.1366595029830933", "0.9761302471160889"),
    @("1.161849286290817E-05", "1", "0.2689805030822754", "0.977703332901001"),
    @("5.964841329841875E-05", "0.9999544024467468", "0.3837735056877136", "0.9662813544273376"),
    @("0.001003601588308811", "0.9998859763145447", "0.1363429874181747", "0.978729248046875"),
    @("0.0004982061218470335", "0.9998631477355957", "0.1847708523273468", "0.978455662727356"),
    @("0.0001415970909874886", "0.9999315738677979", "0.2155129015445709", "0.9766089916229248"),
    @("0.0002008772571571171", "0.9999544024467468", "0.1971027106046677", "0.9778400659561157"),
    @("8.548436017008498E-05", "0.9999544024467468", "0.2450151294469833", "0.9768825769424438"),
    @("1.567982326378115E-05", "1", "0.2663508355617523", "0.9779084920883179"),
    @("0.0009884556056931615", "0.9998403191566467", "0.2534385025501251", "0.9711374044418335"),
    @("0.0002898959792219102", "0.9999544024467468", "0.1425328105688095", "0.9774297475814819"),
    @("0.0002081829588860273", "0.999977171421051", "0.1394363194704056", "0.9764721989631653"),
    @("6.473636858572718E-06", "1", "0.1379337906837463", "0.9771561622619629"),
    @("2.362749910389539E-06", "1", "0.1382739394903183", "0.9775665402412415"),
    @("0.000162509226356633", "0.9999544024467468", "0.09110292792320251", "0.9850899577140808"),
    @("0.0008994809468276799", "0.9998631477355957", "0.1733606457710266", "0.9775665402412415"),
    @("0.0003146666858810931", "0.9999315738677979", "0.1122124791145325", "0.9785924553871155"),
    @("0.0001249395863851532", "0.9999315738677979", "0.1139589548110962", "0.979755163192749"),
    @("1.379447166982573E-05", "1", "0.1144815981388092", "0.9792763590812683"),
    @("0.0003350191691424698", "0.9999315738677979", "0.1558061093091965", "0.9782504439353943"),
    @("3.914447734132409E-05", "0.999977171421051", "0.2084523439407349", "0.9816702008247375"),
    @("0.0001600033865543082", "0.9999315738677979", "0.357175350189209", "0.9692907333374023"),
    @("8.986104148789309E-06", "1", "0.3566190600395203", "0.974078357219696"),
    @("0.0001042160220094956", "0.999977171421051", "0.210339680314064", "0.9789344072341919")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = [double]$data[$i][0]
    $ws.Cells.Item($r, 2).Value = [double]$data[$i][1]
    $ws.Cells.Item($r, 3).Value = [double]$data[$i][2]
    $ws.Cells.Item($r, 4).Value = [double]$data[$i][3]
}